# Add data for 2021-11-07 -- updates "through 10-29" -> "through 10-30"
# for the carjacking-arrests-by-month-yoy workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab to reflect the new "through" date
$ws.Name = "Through 2021-10-30"

# Update row 12 (October row) label
$ws.Range("A12").Value = "October (through 10-30)"

# Row 12 data updates (2018 col block K:M changed; other blocks too)
$ws.Range("C12").Value = 28
$ws.Range("D12").Value = 0.0667

$ws.Range("F12").Value = 50
$ws.Range("G12").Value = 0.1071

$ws.Range("I12").Value = 68
$ws.Range("J12").Value = 0.1392

$ws.Range("K12").Value = 6
$ws.Range("L12").Value = 59
$ws.Range("M12").Value = 0.0923

$ws.Range("N12").Value = 5
$ws.Range("O12").Value = 54
$ws.Range("P12").Value = 0.0847

$ws.Range("R12").Value = 147
$ws.Range("S12").Value = 0.0068

$ws.Range("T12").Value = 2
$ws.Range("U12").Value = 187
$ws.Range("V12").Value = 0.0106

# Row 13 (Total row) data updates
$ws.Range("C13").Value = 224
$ws.Range("D13").Value = 0.125

$ws.Range("F13").Value = 433
$ws.Range("G13").Value = 0.1072

$ws.Range("I13").Value = 645
$ws.Range("J13").Value = 0.0864

$ws.Range("K13").Value = 67
$ws.Range("L13").Value = 546
$ws.Range("M13").Value = 0.1093

$ws.Range("N13").Value = 48
$ws.Range("O13").Value = 433
$ws.Range("P13").Value = 0.0998

$ws.Range("R13").Value = 995
$ws.Range("S13").Value = 0.0515

$ws.Range("T13").Value = 85
$ws.Range("U13").Value = 1352
$ws.Range("V13").Value = 0.0592
